$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.000966519755663806
$ws.Range("C2").Value = 0.0138019021108791
$ws.Range("D2").Value = 0.0871800819608753
$ws.Range("E2").Value = 0.000579911853398283
$ws.Range("F2").Value = 0.971236372071445
$ws.Range("G2").Value = 0.00862135622052115
$ws.Range("H2").Value = 0.994123559885564
$ws.Range("I2").Value = 0.00572179695352973
$ws.Range("J2").Value = 0.00506456351967834
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0.0105543957318488
$ws.Range("M2").Value = 0.0677337044769195
$ws.Range("N2").Value = 0.00100518054589036
$ws.Range("O2").Value = 0.855408644552695
$ws.Range("P2").Value = 0.954109642001083
$ws.Range("Q2").Value = 0.00143044923838243
$ws.Range("R2").Value = 0.0615866388308977
$ws.Range("S2").Value = 0.998337586020258
$ws.Range("T2").Value = 0.998530889971391
$ws.Range("U2").Value = 0.000193303951132761
$ws.Range("V2").Value = 0.178844815588031
$ws.Range("W2").Value = 0.00146911002860898
$ws.Range("X2").Value = 0.89078326760999

# Row 3
$ws.Range("B3").Value = 0.997255083893915
$ws.Range("C3").Value = 0.951712673007036
$ws.Range("D3").Value = 0.910345627464625
$ws.Range("E3").Value = 0.000463929482718627
$ws.Range("F3").Value = 0.00208768267223382
$ws.Range("G3").Value = 0.99114667903812
$ws.Range("H3").Value = 0.000541251063171731
$ws.Range("I3").Value = 0.0000386607902265522
$ws.Range("J3").Value = 0.000541251063171731
$ws.Range("K3").Value = 0.999072141034563
$ws.Range("L3").Value = 0.98766720791773
$ws.Range("M3").Value = 0.930565220753112
$ws.Range("N3").Value = 0.000154643160906209
$ws.Range("O3").Value = 0.00146911002860898
$ws.Range("P3").Value = 0.0443052655996289
$ws.Range("Q3").Value = 0.0000773215804531045
$ws.Range("R3").Value = 0.000154643160906209
$ws.Range("S3").Value = 0.000657233433851388
$ws.Range("T3").Value = 0.000734555014304492
$ws.Range("U3").Value = 0.231810098198407
$ws.Range("V3").Value = 0.819338127271321
$ws.Range("W3").Value = 0.000541251063171731
$ws.Range("X3").Value = 0.0166241397974175

# Row 4
$ws.Range("B4").Value = 0.00139178844815588
$ws.Range("C4").Value = 0.000579911853398283
$ws.Range("D4").Value = 0.000115982370679657
$ws.Range("E4").Value = 0.998646872342071
$ws.Range("F4").Value = 0.0254774607592979
$ws.Range("G4").Value = 0.0000386607902265522
$ws.Range("H4").Value = 0.00502590272945179
$ws.Range("I4").Value = 0.994123559885564
$ws.Range("J4").Value = 0.99427820304647
$ws.Range("K4").Value = 0.0000773215804531045
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0.00069589422407794
$ws.Range("N4").Value = 0.998685533132297
$ws.Range("O4").Value = 0.140454650893064
$ws.Range("P4").Value = 0.00069589422407794
$ws.Range("Q4").Value = 0.998066960488672
$ws.Range("R4").Value = 0.93698291193072
$ws.Range("S4").Value = 0.000541251063171731
$ws.Range("T4").Value = 0.000463929482718627
$ws.Range("U4").Value = 0.000618572643624836
$ws.Range("V4").Value = 0.00034794711203897
$ws.Range("W4").Value = 0.997680352586407
$ws.Range("X4").Value = 0.091432768885796

# Row 5
$ws.Range("B5").Value = 0.000309286321812418
$ws.Range("C5").Value = 0.0330936364339287
$ws.Range("D5").Value = 0.00069589422407794
$ws.Range("E5").Value = 0.000309286321812418
$ws.Range("F5").Value = 0.000889198175210701
$ws.Range("G5").Value = 0.000193303951132761
$ws.Range("H5").Value = 0.000154643160906209
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.0000773215804531045
$ws.Range("K5").Value = 0.000850537384984149
$ws.Range("L5").Value = 0.00154643160906209
$ws.Range("M5").Value = 0.000425268692492075
$ws.Range("N5").Value = 0.000154643160906209
$ws.Range("O5").Value = 0.000309286321812418
$ws.Range("P5").Value = 0.000270625531585866
$ws.Range("Q5").Value = 0.000425268692492075
$ws.Range("R5").Value = 0.0000386607902265522
$ws.Range("S5").Value = 0.000425268692492075
$ws.Range("T5").Value = 0.000270625531585866
$ws.Range("U5").Value = 0.764865073842109
$ws.Range("V5").Value = 0.000463929482718627
$ws.Range("W5").Value = 0.000231964741359313
$ws.Range("X5").Value = 0.0000773215804531045

